# Update view-count-like figures in column F on the "展览" and "全部类型" sheets.
# Both sheets share the same set of events (in the same row order), except
# that "全部类型" has one extra row inserted before the last changed event,
# which shifts that particular row from F24 to F25.

$wb = $excel.ActiveWorkbook

$updates = @{
    "2" = 1580
    "3" = 8935
    "6" = 679
    "8" = 166
    "10" = 64
    "11" = 3795
    "15" = 4157
    "18" = 1138
    "22" = 9
    "23" = 2599
}

# Sheet "展览": last update is on row 24
$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}
$ws1.Range("F24").Value = 102

# Sheet "全部类型": last update is on row 25 (one extra row earlier in the sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates.Keys) {
    $ws4.Range("F$row").Value = $updates[$row]
}
$ws4.Range("F25").Value = 102
